# Insert a new weekly observation row for Orégano (Mercado Mayorista Lo
# Valledor de Santiago) as the new row 138, pushing the existing rows
# 138-240 down to 139-241.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(138).Insert()

$ws.Cells.Item(138, 1).Value = 6
$ws.Cells.Item(138, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(138, 3).Value = "Metropolitana"
$ws.Cells.Item(138, 4).Value = "2022-09-23"
$ws.Cells.Item(138, 5).Value = 13
$ws.Cells.Item(138, 6).Value = 100112029
$ws.Cells.Item(138, 7).Value = "Orégano"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 47
$ws.Cells.Item(138, 11).Value = 15000
$ws.Cells.Item(138, 12).Value = 16000
$ws.Cells.Item(138, 13).Value = 15447
$ws.Cells.Item(138, 14).Value = "`$/docena de atados"
$ws.Cells.Item(138, 15).Value = "Región Metropolitana"
$ws.Cells.Item(138, 16).Value = 5149
$ws.Cells.Item(138, 17).Value = 3
$ws.Cells.Item(138, 18).Value = "Hortaliza"
